$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.461.59'
$ws.Range("E2").Value = '  -3.07%  '
$ws.Range("D3").Value = '3.421.43'
$ws.Range("E3").Value = '  -2.96%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.47'
$ws.Range("E5").Value = '  -3.98%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.35'
$ws.Range("E6").Value = '  -6.36%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '3.422.67'
$ws.Range("E8").Value = '  -2.90%  '
$ws.Range("E9").Value = '  -6.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.120'
$ws.Range("E10").Value = '  -8.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.00'
$ws.Range("E11").Value = '  -9.09%  '
$ws.Range("E12").Value = '  -8.36%  '
$ws.Range("D13").Value = '4.007.52'
$ws.Range("E13").Value = '  -2.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000178'
$ws.Range("E14").Value = '  -8.63%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '3.453.14'
$ws.Range("E15").Value = '  -2.01%  '
$ws.Range("B16").Value = 'TRON'
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.115'
$ws.Range("E16").Value = '  -1.50%  '
$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.19'
$ws.Range("E17").Value = '  -8.66%  '
$ws.Range("D18").Value = '64.454.11'
$ws.Range("E18").Value = '  -2.89%  '
$ws.Range("E19").Value = '  -11.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.67'
$ws.Range("E20").Value = '  -8.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.54'
$ws.Range("E21").Value = '  -7.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '379.48'
$ws.Range("E22").Value = '  -10.39%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.543'
$ws.Range("E23").Value = '  -8.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.73'
$ws.Range("E25").Value = '  -0.22%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '71.76'
$ws.Range("E26").Value = '  -7.11%  '
$ws.Range("D27").Value = '3.563.96'
$ws.Range("E27").Value = '  -2.88%  '
$ws.Range("E28").Value = '  -8.38%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.14'
$ws.Range("E30").Value = '  -9.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.04'
$ws.Range("E31").Value = '  -10.14%  '
$ws.Range("E32").Value = '  -11.23%  '
$ws.Range("D33").Value = '3.439.42'
$ws.Range("E33").Value = '  -2.70%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '22.95'
$ws.Range("E35").Value = '  -5.28%  '
$ws.Range("E36").Value = '  -9.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '171.12'
$ws.Range("E37").Value = '  -1.87%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.17'
$ws.Range("E38").Value = '  -12.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.69'
$ws.Range("E39").Value = '  -11.23%  '
$ws.Range("E40").Value = '  -10.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.65'
$ws.Range("E41").Value = '  -10.94%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0757'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.798'
$ws.Range("E43").Value = '  -6.95%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.53'
$ws.Range("E45").Value = '  -8.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.27'
$ws.Range("E46").Value = '  -14.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.58'
$ws.Range("E47").Value = '  -10.26%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.09'
$ws.Range("E48").Value = '  -1.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.46'
$ws.Range("E49").Value = '  -1.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.50'
$ws.Range("E50").Value = '  -8.05%  '
$ws.Range("D51").Value = '2.194.63'
$ws.Range("E51").Value = '  -5.17%  '
